$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "requisite stream depth" column (K) previously just rounded J/H up to the
# next whole multiple-count (CEILING(J/H,1)); now it also multiplies back by H
# so the result reflects the theoretical forays in actual stream-depth units
# instead of the raw multiple count.
$ws.Range("K2:K121").Formula = "=CEILING(J2/H2,1)*H2"

# Bring the new wins near the bottom of the sheet into view and leave the
# selection where the analysis now continues.
$ws.Activate()
$ws.Range("K124").Select()
